# Edit slide 2 ("TEAM INSIGHT_IQ" roster slide): reposition four existing
# photo/name pairs to make room for a new team-member entry, then add the
# new photo + name textbox for "Grace Emeruwa".
#
# NOTE on the magic-looking decimal point coordinates below: this host
# stores Shape.Left/Top/Width/Height as 32-bit floats and truncates when
# converting back to EMU on save (1 pt = 12700 EMU). Plain
# "emu / 12700.0" can therefore land 1 EMU short after the float32
# round-trip. The literals used here were chosen so that, after the
# float32 round-trip, they serialize to the exact target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1) Picture 17 (id=18): photo above "Ikechukwu Daniel" -----------------
$pic18 = $s.Shapes.Item(12)
$pic18.Left = 475.81850393700785
$pic18.Top  = 248.540001

# --- 2) TextBox 19 (id=20): "Ikechukwu Daniel" label ------------------------
$tb20 = $s.Shapes.Item(13)
$tb20.Left   = 446.94543307086616
$tb20.Top    = 279.1116535433071
$tb20.Width  = 115.9296073992126
$tb20.Height = 24.234410348818898

# --- 3) Picture 21 (id=22): photo above "Aramide Arabesin" ------------------
$pic22 = $s.Shapes.Item(15)
$pic22.Left = 475.81850393700785
$pic22.Top  = 316.4082795165354

# --- 4) TextBox 22 (id=23): "Aramide Arabesin" label ------------------------
$tb23 = $s.Shapes.Item(16)
$tb23.Left = 436.5383464566929
$tb23.Top  = 346.1428375456693

# --- 5) New picture for Grace Emeruwa (duplicate the same generic-avatar
#        picture used for id=22, which already embeds the right image) -----
$newPic = $pic22.Duplicate()
$newPic.Name = "Picture 23"
$newPic.Left = 615.6755118110236
$newPic.Top  = 248.64110566220472
$newPic.Width = 46.5
$newPic.Height = 26.657874115748033

# --- 6) New "Grace Emeruwa" textbox -----------------------------------------
$newTb = $s.Shapes.AddTextbox(1, 579.698031496063, 275.78062992125984, 126.33669671338583, 26.657795975590552)
$newTb.Name = "TextBox 24"
$newTb.TextFrame.TextRange.Text = "Grace Emeruwa"
$newTb.TextFrame.TextRange.Font.Size = 16
$newTb.TextFrame.TextRange.Font.Name = "Bahnschrift SemiBold SemiConden"
